$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update A3 input value (0.28 instead of 10.65) ---
# This ripples through B3 (=A3*A4/(A4+A5)) and C3 (=ROUND(B3*1024/A6,0))
$ws.Range("A3").Value = 0.28000000000000003

# --- New row 8: a second voltage-divider style pickup used for battery level ---
$ws.Range("A8").Value = 0.5
$ws.Range("B8").Formula = "=A3/A8"

# Match A8's look (yellow fill / General number format) to the other input
# cells in column A (A3:A6) by copying formats from A6.
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Leave the active cell / selection on A3 like the saved workbook ---
$ws.Range("A3").Select()
